# The edit reorders data rows 3-14 (columns A:AY) of the active sheet.
# Row contents themselves are unchanged - only their row positions move,
# per this after-row -> before-row mapping (1-based worksheet rows):
#   3<-7  4<-8  5<-9  6<-3  7<-10  8<-4  9<-5  10<-11  11<-12  12<-13  13<-14  14<-6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow = 14
$lastCol = "AY"

# Snapshot every source row (A:AY) into memory before writing anything,
# so overlapping reads/writes in the cyclic permutation don't clobber data.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rng = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $snapshot[$r] = $rng.Value()
}

# after-row -> before-row mapping
$mapping = @{
    3  = 7
    4  = 8
    5  = 9
    6  = 3
    7  = 10
    8  = 4
    9  = 5
    10 = 11
    11 = 12
    12 = 13
    13 = 14
    14 = 6
}

# Columns that hold values which LOOK numeric/date-like but are stored as
# plain text in the workbook (e.g. "2021-04-25", "00:00", "1"). Excel's
# Range.Value setter auto-detects and silently converts such text into real
# date/number values, which would corrupt the data on a round trip. Force
# those destination columns to Text format first so the setter leaves the
# strings alone.
$textProtectCols = @("I", "Y", "Z", "AA", "AB")

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapping[$r]
    foreach ($col in $textProtectCols) {
        $ws.Range($col + $r).NumberFormat = "@"
    }
    $destRng = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $destRng.Value = $snapshot[$srcRow]
}
